$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'319.79"
$ws.Range("E2").Value = "'-3.66%"
$ws.Range("D3").Value = "'42.59"
$ws.Range("E3").Value = "'-7.11%"
$ws.Range("D4").Value = "'5.234"
$ws.Range("E4").Value = "'-7.67%"
$ws.Range("D5").Value = "'0.08162"
$ws.Range("E5").Value = "'-2.48%"
$ws.Range("D6").Value = "'4.329"
$ws.Range("E6").Value = "'-3.36%"
$ws.Range("D7").Value = "'1.812"
$ws.Range("E7").Value = "'-11.14%"
$ws.Range("D8").Value = "'0.9515"
$ws.Range("E8").Value = "'-3.43%"
$ws.Range("D9").Value = "'0.1111"
$ws.Range("E9").Value = "'-3.90%"
$ws.Range("D10").Value = "'0.1848"
$ws.Range("E10").Value = "'-4.77%"
$ws.Range("D11").Value = "'0.09341"
$ws.Range("E11").Value = "'-7.50%"
$ws.Range("D12").Value = "'0.04656"
$ws.Range("E12").Value = "'-0.30%"
$ws.Range("D13").Value = "'7.447"
$ws.Range("E13").Value = "'-28.02%"
$ws.Range("D14").Value = "'0.1060"
$ws.Range("E14").Value = "'0.15%"
$ws.Range("D15").Value = "'0.001281"
$ws.Range("E15").Value = "'0.15%"
$ws.Range("D16").Value = "'0.005819"
$ws.Range("E16").Value = "'-3.46%"
$ws.Range("E17").Value = "'-0.05%"
$ws.Range("D18").Value = "'2.517"
$ws.Range("E18").Value = "'-2.66%"
$ws.Range("E19").Value = "'-0.06%"
$ws.Range("D20").Value = "'0.1391"
$ws.Range("E20").Value = "'-0.62%"
$ws.Range("D21").Value = "'0.2686"
$ws.Range("E21").Value = "'3.07%"
$ws.Range("D22").Value = "'0.04197"
$ws.Range("E23").Value = "'-4.17%"
$ws.Range("D24").Value = "'0.004318"
$ws.Range("E24").Value = "'-8.24%"
$ws.Range("E25").Value = "'1.89%"
$ws.Range("E26").Value = "'-20.11%"
$ws.Range("D38").Value = "'0.02593"
$ws.Range("E38").Value = "'-7.10%"
$ws.Range("D39").Value = "'0.05478"
$ws.Range("E39").Value = "'-5.47%"
$ws.Range("D40").Value = "'0.007780"
$ws.Range("E40").Value = "'0.52%"
$ws.Range("D41").Value = "'0.1394"
$ws.Range("E41").Value = "'-2.99%"
$ws.Range("D42").Value = "'0.006626"
$ws.Range("E42").Value = "'-7.86%"
$ws.Range("E43").Value = "'7.81%"
$ws.Range("D44").Value = "'0.008489"
$ws.Range("E44").Value = "'4.99%"
$ws.Range("D45").Value = "'0.3425"
$ws.Range("E45").Value = "'-2.16%"
$ws.Range("D46").Value = "'0.00006981"
$ws.Range("E46").Value = "'-5.64%"
$ws.Range("E47").Value = "'0.30%"
$ws.Range("D48").Value = "'0.003466"
$ws.Range("E48").Value = "'-1.08%"
$ws.Range("D49").Value = "'0.003544"
$ws.Range("E49").Value = "'1.22%"
$ws.Range("E50").Value = "'0.30%"
$ws.Range("E51").Value = "'0.30%"
